# Weekly fruit/vegetable price update: insert the latest week's
# observations (2023-05-30) for "Terminal Hortofrutícola Agro Chillán -
# Membrillo" at the top of the data block (row 6), pushing the existing
# history down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 6, shifting rows 6:24 down to 8:26
# (xlShiftDown = -4121)
$ws.Range("A6:T7").Insert(-4121)

# Row 6: Membrillo Champion, Primera
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 45076
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100104
$ws.Range("H6").Value = "Frutos de pepita"
$ws.Range("I6").Value = 100104003
$ws.Range("J6").Value = "Membrillo"
$ws.Range("K6").Value = "Champion"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = "$/caja 15 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 800
$ws.Range("T6").Value = 15

# Row 7: Membrillo Champion, Segunda
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 45076
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100104
$ws.Range("H7").Value = "Frutos de pepita"
$ws.Range("I7").Value = 100104003
$ws.Range("J7").Value = "Membrillo"
$ws.Range("K7").Value = "Champion"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 10000
$ws.Range("Q7").Value = "$/caja 15 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 667
$ws.Range("T7").Value = 15
